$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 128 ("Fruta / hortaliza, semanal"),
# pushing the existing rows 128-155 down to 129-156.
$ws.Rows("128:128").Insert()

$ws.Cells.Item(128, 1).Value  = 8
$ws.Cells.Item(128, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(128, 3).Value  = "Coquimbo"
$ws.Cells.Item(128, 4).Value  = 44551
$ws.Cells.Item(128, 5).Value  = 4
$ws.Cells.Item(128, 6).Value  = 100112021
$ws.Cells.Item(128, 7).Value  = "Ají"
$ws.Cells.Item(128, 8).Value  = "Inferno"
$ws.Cells.Item(128, 9).Value  = "Primera"
$ws.Cells.Item(128, 10).Value = 500
$ws.Cells.Item(128, 11).Value = 15500
$ws.Cells.Item(128, 12).Value = 16000
$ws.Cells.Item(128, 13).Value = 15750
$ws.Cells.Item(128, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 1050
$ws.Cells.Item(128, 17).Value = 15
$ws.Cells.Item(128, 18).Value = "Hortaliza"
